$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = '@'
    $Range.Value = $Text
    $Range.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '63.719.37'
Set-TextValue $ws.Range('E2') '  -1.26%  '
Set-TextValue $ws.Range('D3') '2.638.04'
Set-TextValue $ws.Range('E3') '  +0.17%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '579.09'
Set-TextValue $ws.Range('E5') '  -0.03%  '
Set-TextValue $ws.Range('D6') '155.20'
Set-TextValue $ws.Range('E6') '  -0.98%  '
Set-TextValue $ws.Range('E7') '  +0.06%  '
Set-TextValue $ws.Range('D8') '0.620'
Set-TextValue $ws.Range('E8') '  -4.27%  '
Set-TextValue $ws.Range('D9') '2.636.33'
Set-TextValue $ws.Range('E9') '  +0.22%  '
Set-TextValue $ws.Range('E10') '  -4.16%  '
Set-TextValue $ws.Range('D11') '5.80'
Set-TextValue $ws.Range('E11') '  -0.14%  '
Set-TextValue $ws.Range('E12') '  -2.22%  '
Set-TextValue $ws.Range('E13') '  +0.91%  '
Set-TextValue $ws.Range('D14') '28.39'
Set-TextValue $ws.Range('E14') '  -0.84%  '
Set-TextValue $ws.Range('D15') '3.115.62'
Set-TextValue $ws.Range('E15') '  +0.23%  '
Set-TextValue $ws.Range('E16') '  -2.64%  '
Set-TextValue $ws.Range('D17') '63.703.11'
Set-TextValue $ws.Range('E17') '  -0.91%  '
Set-TextValue $ws.Range('D18') '2.648.13'
Set-TextValue $ws.Range('E18') '  +0.27%  '
Set-TextValue $ws.Range('E19') '  -1.32%  '
Set-TextValue $ws.Range('E20') '  +3.02%  '
Set-TextValue $ws.Range('E21') '  -3.36%  '
Set-TextValue $ws.Range('D22') '344.29'
Set-TextValue $ws.Range('E22') '  -0.87%  '
Set-TextValue $ws.Range('E23') '  +0.28%  '
Set-TextValue $ws.Range('D24') '68.02'
Set-TextValue $ws.Range('E24') '  +0.32%  '
Set-TextValue $ws.Range('D25') '1.87'
Set-TextValue $ws.Range('E25') '  +7.32%  '
Set-TextValue $ws.Range('E26') '  -4.88%  '
Set-TextValue $ws.Range('D27') '602.20'
Set-TextValue $ws.Range('E27') '  +5.98%  '
Set-TextValue $ws.Range('D28') '9.23'
Set-TextValue $ws.Range('E28') '  -1.79%  '
Set-TextValue $ws.Range('D29') '1.60'
Set-TextValue $ws.Range('E29') '  +1.62%  '
Set-TextValue $ws.Range('D30') '8.17'
Set-TextValue $ws.Range('E30') '  +3.00%  '
Set-TextValue $ws.Range('E31') '  -0.59%  '
Set-TextValue $ws.Range('E32') '  +0.13%  '
Set-TextValue $ws.Range('E33') '  -0.50%  '
Set-TextValue $ws.Range('E34') '  +0.84%  '
Set-TextValue $ws.Range('D35') '6.56'
Set-TextValue $ws.Range('E35') '  -2.16%  '
Set-TextValue $ws.Range('D36') '5.43'
Set-TextValue $ws.Range('E36') '  +2.24%  '
Set-TextValue $ws.Range('E37') '  -2.51%  '
Set-TextValue $ws.Range('D38') '1.00'
Set-TextValue $ws.Range('E38') '  +0.10%  '
Set-TextValue $ws.Range('D39') '19.71'
Set-TextValue $ws.Range('E39') '  -1.95%  '
Set-TextValue $ws.Range('D40') '1.89'
Set-TextValue $ws.Range('E40') '  -2.42%  '
Set-TextValue $ws.Range('D41') '151.02'
Set-TextValue $ws.Range('E41') '  -2.60%  '
Set-TextValue $ws.Range('B42') 'USDe'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D42') '0.999'
Set-TextValue $ws.Range('E42') '  -0.03%  '
Set-TextValue $ws.Range('B43') 'dogwifhat'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D43') '2.55'
Set-TextValue $ws.Range('E43') '  +3.33%  '
Set-TextValue $ws.Range('E44') '  -0.73%  '
Set-TextValue $ws.Range('D45') '160.59'
Set-TextValue $ws.Range('E45') '  +1.26%  '
Set-TextValue $ws.Range('D46') '24.15'
Set-TextValue $ws.Range('E46') '  +4.97%  '
Set-TextValue $ws.Range('E47') '  -2.30%  '
Set-TextValue $ws.Range('D48') '0.0585'
Set-TextValue $ws.Range('E48') '  -2.60%  '
Set-TextValue $ws.Range('D49') '0.632'
Set-TextValue $ws.Range('E49') '  -0.68%  '
Set-TextValue $ws.Range('B50') 'Stellar'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D50') '0.0995'
Set-TextValue $ws.Range('E50') '  -2.65%  '
Set-TextValue $ws.Range('B51') 'VeChain'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D51') '0.0248'
Set-TextValue $ws.Range('E51') '  -1.37%  '
